$p = $ppt.ActivePresentation
$s = $p.Slides.Item(13)

# --- Shape 1123 ("Rect 0") : nudge position/size, merge adjacent same-format runs ---
$sh = $s.Shapes.Item(3)

# Position/size tweak (values picked so the EMU round-trip lands on the exact target)
$sh.Left   = 536.10001
$sh.Width  = 325.40001
$sh.Height = 75.25

# Merge runs that share identical formatting (pure run-splitting cleanup, text unchanged)
$tr = $sh.TextFrame.TextRange
$tr.Characters(1, 2).Text  = "31"
$tr.Characters(34, 3).Text = "함수를"
$tr.Characters(37, 16).Text = " 호출하고 삭제할 오브젝트와 "
$tr.Characters(56, 11).Text = " 시간을 설정합니다."

# --- Drop the stale local-file "descr" (AlternativeText) from the three pictures ---
$s.Shapes.Item(4).AlternativeText = ""
$s.Shapes.Item(5).AlternativeText = ""
$s.Shapes.Item(7).AlternativeText = ""
